# Fix typos in header labels and update the active selection.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "Greater than 5"
$ws.Range("D1").Value = "Divided by 5"

$ws.Range("B1").Select()
